{"js": "// Replace the date and all two-digit multiplication answers in the table.\nconst replacements = [\n  [\"2023-10-22 Sunday\", \"2023-10-23 Monday\"],\n  [\"68\\u00d782=5576\", \"69\\u00d746=3174\"],\n  [\"50\\u00d774=3700\", \"24\\u00d763=1512\"],\n  [\"64\\u00d742=2688\", \"99\\u00d720=1980\"],\n  [\"34\\u00d751=1734\", \"24\\u00d795=2280\"],\n  [\"51\\u00d791=4641\", \"70\\u00d731=2170\"],\n  [\"76\\u00d727=2052\", \"39\\u00d740=1560\"],\n  [\"83\\u00d772=5976\", \"56\\u00d747=2632\"],\n  [\"41\\u00d714=574\", \"97\\u00d781=7857\"],\n  [\"17\\u00d735=595\", \"14\\u00d739=546\"],\n  [\"16\\u00d746=736\", \"51\\u00d723=1173\"],\n  [\"21\\u00d741=861\", \"55\\u00d769=3795\"],\n  [\"25\\u00d776=1900\", \"32\\u00d712=384\"],\n  [\"87\\u00d769=6003\", \"81\\u00d777=6237\"],\n  [\"71\\u00d760=4260\", \"26\\u00d795=2470\"],\n  [\"81\\u00d775=6075\", \"74\\u00d756=4144\"],\n  [\"72\\u00d785=6120\", \"82\\u00d750=4100\"],\n  [\"36\\u00d750=1800\", \"19\\u00d740=760\"],\n  [\"26\\u00d732=832\", \"38\\u00d772=2736\"],\n  [\"49\\u00d767=3283\", \"52\\u00d789=4628\"],\n  [\"39\\u00d767=2613\", \"27\\u00d772=1944\"],\n  [\"53\\u00d769=3657\", \"65\\u00d749=3185\"],\n  [\"13\\u00d718=234\", \"96\\u00d756=5376\"],\n  [\"55\\u00d752=2860\", \"66\\u00d744=2904\"],\n  [\"55\\u00d796=5280\", \"64\\u00d740=2560\"],\n  [\"56\\u00d789=4984\", \"23\\u00d789=2047\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-10-22 Sunday\", \"2023-10-23 Monday\"),\n    @(\"68\u00d782=5576\", \"69\u00d746=3174\"),\n    @(\"50\u00d774=3700\", \"24\u00d763=1512\"),\n    @(\"64\u00d742=2688\", \"99\u00d720=1980\"),\n    @(\"34\u00d751=1734\", \"24\u00d795=2280\"),\n    @(\"51\u00d791=4641\", \"70\u00d731=2170\"),\n    @(\"76\u00d727=2052\", \"39\u00d740=1560\"),\n    @(\"83\u00d772=5976\", \"56\u00d747=2632\"),\n    @(\"41\u00d714=574\", \"97\u00d781=7857\"),\n    @(\"17\u00d735=595\", \"14\u00d739=546\"),\n    @(\"16\u00d746=736\", \"51\u00d723=1173\"),\n    @(\"21\u00d741=861\", \"55\u00d769=3795\"),\n    @(\"25\u00d776=1900\", \"32\u00d712=384\"),\n    @(\"87\u00d769=6003\", \"81\u00d777=6237\"),\n    @(\"71\u00d760=4260\", \"26\u00d795=2470\"),\n    @(\"81\u00d775=6075\", \"74\u00d756=4144\"),\n    @(\"72\u00d785=6120\", \"82\u00d750=4100\"),\n    @(\"36\u00d750=1800\", \"19\u00d740=760\"),\n    @(\"26\u00d732=832\", \"38\u00d772=2736\"),\n    @(\"49\u00d767=3283\", \"52\u00d789=4628\"),\n    @(\"39\u00d767=2613\", \"27\u00d772=1944\"),\n    @(\"53\u00d769=3657\", \"65\u00d749=3185\"),\n    @(\"13\u00d718=234\", \"96\u00d756=5376\"),\n    @(\"55\u00d752=2860\", \"66\u00d744=2904\"),\n    @(\"55\u00d796=5280\", \"64\u00d740=2560\"),\n    @(\"56\u00d789=4984\", \"23\u00d789=2047\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Execute(\n        $find,       # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $replace,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
